# Updates the cryptos list (Price / Volume(1h) columns) on the active sheet.
# Values are taken from the authoritative diff describing this commit.
# Some "Price" values look like plain numbers (e.g. "30.61"); Excel would
# otherwise silently convert them to numeric cells, so those specific cells
# are forced to Text format before the value is written so they stay
# strings, matching the source data (e.g. "67.110.14" stays text naturally
# because it has more than one '.' and can never parse as a number).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = "67.110.14"; E = "  -0.26%  "; ForceText = $false },
    @{ Row = 3; D = "3.535.56"; E = "  +1.50%  "; ForceText = $false },
    @{ Row = 4; D = $null; E = "  -0.04%  "; ForceText = $false },
    @{ Row = 5; D = $null; E = "  -0.77%  "; ForceText = $false },
    @{ Row = 6; D = "177.94"; E = "  -0.16%  "; ForceText = $true },
    @{ Row = 7; D = $null; E = "  +0.00%  "; ForceText = $false },
    @{ Row = 8; D = "3.533.47"; E = "  +1.35%  "; ForceText = $false },
    @{ Row = 9; D = $null; E = "  +0.28%  "; ForceText = $false },
    @{ Row = 10; D = $null; E = "  -0.97%  "; ForceText = $false },
    @{ Row = 11; D = $null; E = "  -1.53%  "; ForceText = $false },
    @{ Row = 12; D = $null; E = "  -2.01%  "; ForceText = $false },
    @{ Row = 13; D = "4.147.09"; E = "  +1.50%  "; ForceText = $false },
    @{ Row = 14; D = "30.61"; E = "  -4.01%  "; ForceText = $true },
    @{ Row = 15; D = $null; E = "  -2.26%  "; ForceText = $false },
    @{ Row = 16; D = "67.066.85"; E = "  -0.35%  "; ForceText = $false },
    @{ Row = 17; D = $null; E = "  -1.17%  "; ForceText = $false },
    @{ Row = 18; D = "3.531.95"; E = "  +1.44%  "; ForceText = $false },
    @{ Row = 19; D = $null; E = "  -1.86%  "; ForceText = $false },
    @{ Row = 20; D = $null; E = "  -1.16%  "; ForceText = $false },
    @{ Row = 21; D = "384.85"; E = "  -0.91%  "; ForceText = $true },
    @{ Row = 22; D = $null; E = "  -1.18%  "; ForceText = $false },
    @{ Row = 23; D = "0.544"; E = "  +1.34%  "; ForceText = $true },
    @{ Row = 24; D = $null; E = "  +0.10%  "; ForceText = $false },
    @{ Row = 25; D = $null; E = "  +0.73%  "; ForceText = $false },
    @{ Row = 26; D = "72.04"; E = "  -2.44%  "; ForceText = $true },
    @{ Row = 27; D = $null; E = "  +1.62%  "; ForceText = $false },
    @{ Row = 28; D = "9.99"; E = "  -3.35%  "; ForceText = $true },
    @{ Row = 29; D = $null; E = "  +0.17%  "; ForceText = $false },
    @{ Row = 30; D = $null; E = "  +0.02%  "; ForceText = $false },
    @{ Row = 31; D = "24.61"; E = "  +4.70%  "; ForceText = $true },
    @{ Row = 32; D = "5.98"; E = "  -2.21%  "; ForceText = $true },
    @{ Row = 33; D = $null; E = "  -1.16%  "; ForceText = $false },
    @{ Row = 34; D = $null; E = "  -3.49%  "; ForceText = $false },
    @{ Row = 35; D = $null; E = "  -0.87%  "; ForceText = $false },
    @{ Row = 36; D = $null; E = "  -0.08%  "; ForceText = $false },
    @{ Row = 37; D = $null; E = "  -0.05%  "; ForceText = $false },
    @{ Row = 38; D = "29.74"; E = "  +13.84%  "; ForceText = $true },
    @{ Row = 39; D = "159.89"; E = "  -2.87%  "; ForceText = $true },
    @{ Row = 40; D = "0.898"; E = "  +3.39%  "; ForceText = $true },
    @{ Row = 41; D = $null; E = "  -2.58%  "; ForceText = $false },
    @{ Row = 42; D = "6.65"; E = "  -2.39%  "; ForceText = $true },
    @{ Row = 43; D = $null; E = "  -5.14%  "; ForceText = $false },
    @{ Row = 44; D = "4.56"; E = "  -2.08%  "; ForceText = $true },
    @{ Row = 45; D = "2.757.59"; E = "  -2.66%  "; ForceText = $false },
    @{ Row = 46; D = $null; E = "  -1.39%  "; ForceText = $false },
    @{ Row = 47; D = "25.66"; E = "  -5.04%  "; ForceText = $true },
    @{ Row = 48; D = "40.81"; E = "  -1.90%  "; ForceText = $true },
    @{ Row = 49; D = $null; E = "  -0.03%  "; ForceText = $false },
    @{ Row = 50; D = "329.29"; E = "  -1.63%  "; ForceText = $true },
    @{ Row = 51; D = $null; E = "  -1.66%  "; ForceText = $false }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $cell = $ws.Cells.Item($u.Row, 4)
        if ($u.ForceText) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $u.D
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
